# Generate Report for Handoff
#
# Moves the localization status from "In Translation" to "Ready for
# handoff" and refreshes the handoff timestamps, on all three sheets
# (Overview, zh-cn, de-de). Also widens the Status-related columns that
# now need to fit the longer "Ready for handoff" label.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" -----------------
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# --- Timestamps --------------------------------------------------------
# Overview's "Latest HO Xliff Generate Date" and de-de's "Latest Handoff
# Datetime" shared the same instant.
$wsOverview.Range("G2").Value = "2016-09-02 08:45:42"
$wsDeDe.Range("H2").Value     = "2016-09-02 08:45:42"

# zh-cn's "Latest Handoff Datetime" advanced separately.
$wsZhCn.Range("H2").Value = "2016-09-02 08:45:37"

# --- Column widths -------------------------------------------------
# The Status columns widen to fit "Ready for handoff". The host only
# supports whole-pixel column widths (1/6-character granularity), so
# 16.33 is the nearest settable value that lands on the same rounded
# width Excel itself would have produced here.
$wsOverview.Range("E1").ColumnWidth = 16.33
$wsOverview.Range("F1").ColumnWidth = 16.33
$wsZhCn.Range("C1").ColumnWidth     = 16.33
$wsDeDe.Range("C1").ColumnWidth     = 16.33
